# Toggle the "execute" flag for the loginLogoutTest and amazonTest rows,
# and update the last selected cell on each sheet.

$wb = $excel.ActiveWorkbook

# --- RUNMANAGER sheet ---
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws1.Range("C2").Value = "no"    # loginLogoutTest execute: yes -> no
$ws1.Range("C4").Value = "yes"   # amazonTest execute: no -> yes
$ws1.Range("C4").Select()

# --- DATA sheet ---
$ws2 = $wb.Worksheets.Item("DATA")
$ws2.Range("B2").Value = "no"    # loginLogoutTest execute: yes -> no
$ws2.Range("B3").Value = "no"    # loginLogoutTest execute: yes -> no
$ws2.Range("B7").Value = "yes"   # amazonTest execute: no -> yes
$ws2.Range("B4").Select()
